$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "Piège"
$ws.Range("B2").Value = "Comment déjouer le piège"

# Row 4
$ws.Range("A4").Value = "Rivière"
$ws.Range("B4").Value = "Création / déplacement d'un pont"

# Row 5
$ws.Range("A5").Value = "Racine au sol"
$ws.Range("B5").Value = "Destruction des racines"

# Row 6
$ws.Range("A6").Value = "Feuille grimpante"
$ws.Range("B6").Value = "Découper"

# Row 7
$ws.Range("A7").Value = "Rocher au sol"
$ws.Range("B7").Value = "Détruire le rocher"

# Row 8
$ws.Range("A8").Value = "Singe qui lance des noix de coco"
$ws.Range("B8").Value = "Placer un bouclier ou autre pour faire ""Parapluie"""

# Row 9
$ws.Range("A9").Value = "Boulet qui tombe d'une montagne"
$ws.Range("B9").Value = "Cliquer pour le réduire / détruire"

# Row 10
$ws.Range("A10").Value = "Champignon géant -> mur"
$ws.Range("B10").Value = "Cliquer sur le champi pour projeter le perso"

# Row 11
$ws.Range("A11").Value = "Boue sable mouvant"
$ws.Range("B11").Value = "Lui donner une corde pour l'aider"

# Row 12
$ws.Range("A12").Value = "Feu "
$ws.Range("B12").Value = "Cliquer pour réduire la vie du feu"

# Row 13
$ws.Range("A13").Value = "Attaque distance"
$ws.Range("B13").Value = "Cliquer pour détruire/tuer"

# Row 14
$ws.Range("A14").Value = "Attaque à bout portant "
$ws.Range("B14").Value = """"""

# Row 15
$ws.Range("A15").Value = "Attaque venant du ciel (pluie, bulles, cercles de fumée, animaux, objets, nourriture...)"
$ws.Range("B15").Value = "Cliquer pour enlever"

# Row 16
$ws.Range("A16").Value = "Qui foncent en ligne droite ou zigzag sur le personnage ("""")"
$ws.Range("B16").Value = """"""

# Row 17
$ws.Range("A17").Value = "Objets roulants : Boule de neige (montagne), tronc d'arbre (foret)"
$ws.Range("B17").Value = "Cliquer pour dévier les objets"

# Row 18 - B18 no longer used
$ws.Range("A18").Value = "Sol modifié  : Boue -> pieds qui s'enfoncent (foret), marshmallow nuage (sol collant caramel), glisse (montagne)"
$ws.Range("B18").ClearContents()

# Row 19
$ws.Range("A19").Value = "Vide/trou"
$ws.Range("B19").Value = "Mettre une planche"

# Update selection from B18 to A18
$ws.Range("A18").Select()
